# Update Facebook and Twitter pivot table statistics (msg_count_twitter,
# msg_count_twitter_engage, msg_count_facebook blocks) for rows 4-14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("AC4").Value = 55
$ws.Range("AD4").Value = 55
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 55
$ws.Range("AH4").Value = 55
$ws.Range("AI4").Value = 55
$ws.Range("AJ4").Value = 55
$ws.Range("AL4").Value = 55
$ws.Range("AM4").Value = 1
$ws.Range("AN4").Value = 100
$ws.Range("AO4").Value = 2.1
# Row 5
$ws.Range("C5").Value = 1470357
$ws.Range("D5").Value = 1170.7
$ws.Range("E5").Value = 2595.5
$ws.Range("H5").Value = 398
$ws.Range("I5").Value = 1372.8
$ws.Range("J5").Value = 55081
$ws.Range("L5").Value = 1610.5
$ws.Range("M5").Value = 913
$ws.Range("N5").Value = 72.7
$ws.Range("P5").Value = 541002
$ws.Range("Q5").Value = 430.7
$ws.Range("R5").Value = 1830.2
$ws.Range("U5").Value = 64
$ws.Range("V5").Value = 297.5
$ws.Range("W5").Value = 48717
$ws.Range("Y5").Value = 583.6
$ws.Range("Z5").Value = 927
$ws.Range("AA5").Value = 73.8
$ws.Range("AB5").Value = 1.5
$ws.Range("AC5").Value = 682778
$ws.Range("AD5").Value = 543.6
$ws.Range("AE5").Value = 668.9
$ws.Range("AH5").Value = 328.5
$ws.Range("AI5").Value = 856.8
$ws.Range("AJ5").Value = 5487
$ws.Range("AL5").Value = 738.1
$ws.Range("AM5").Value = 925
$ws.Range("AN5").Value = 73.59999999999999
$ws.Range("AO5").Value = 0.7
# Row 6
$ws.Range("C6").Value = 36873
$ws.Range("D6").Value = 899.3
$ws.Range("E6").Value = 3581.2
$ws.Range("I6").Value = 675
$ws.Range("J6").Value = 22972
$ws.Range("L6").Value = 2048.5
$ws.Range("M6").Value = 18
$ws.Range("N6").Value = 43.9
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 20377
$ws.Range("Q6").Value = 497
$ws.Range("R6").Value = 2757.8
$ws.Range("V6").Value = 80
$ws.Range("W6").Value = 17700
$ws.Range("Y6").Value = 1072.5
$ws.Range("Z6").Value = 19
$ws.Range("AA6").Value = 46.3
$ws.Range("AB6").Value = 0.1
$ws.Range("AC6").Value = 14135
$ws.Range("AD6").Value = 344.8
$ws.Range("AE6").Value = 425.6
$ws.Range("AH6").Value = 124
$ws.Range("AI6").Value = 717
$ws.Range("AJ6").Value = 1516
$ws.Range("AL6").Value = 589
$ws.Range("AM6").Value = 24
$ws.Range("AN6").Value = 58.5
$ws.Range("AO6").Value = 0
# Row 7
$ws.Range("C7").Value = 269530
$ws.Range("D7").Value = 1069.6
$ws.Range("E7").Value = 3569.8
$ws.Range("H7").Value = 62.5
$ws.Range("I7").Value = 1033.8
$ws.Range("J7").Value = 50844
$ws.Range("L7").Value = 1911.6
$ws.Range("M7").Value = 141
$ws.Range("N7").Value = 56
$ws.Range("O7").Value = 0.6
$ws.Range("P7").Value = 84346
$ws.Range("Q7").Value = 334.7
$ws.Range("R7").Value = 1284.3
$ws.Range("U7").Value = 9.5
$ws.Range("V7").Value = 219.2
$ws.Range("W7").Value = 11140
$ws.Range("Y7").Value = 594
$ws.Range("Z7").Value = 142
$ws.Range("AA7").Value = 56.3
$ws.Range("AC7").Value = 123602
$ws.Range("AD7").Value = 490.5
$ws.Range("AE7").Value = 857.8
$ws.Range("AH7").Value = 130
$ws.Range("AI7").Value = 703.8
$ws.Range("AJ7").Value = 8295
$ws.Range("AL7").Value = 763
$ws.Range("AM7").Value = 162
$ws.Range("AN7").Value = 64.3
$ws.Range("AO7").Value = 0.3
# Row 8
$ws.Range("C8").Value = 72063
$ws.Range("D8").Value = 713.5
$ws.Range("E8").Value = 1214.9
$ws.Range("H8").Value = 166
$ws.Range("I8").Value = 948
$ws.Range("J8").Value = 6588
$ws.Range("L8").Value = 1201
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 59.4
$ws.Range("P8").Value = 13377
$ws.Range("Q8").Value = 132.4
$ws.Range("R8").Value = 268.3
$ws.Range("U8").Value = 15
$ws.Range("V8").Value = 158
$ws.Range("W8").Value = 1567
$ws.Range("Y8").Value = 219.3
$ws.Range("Z8").Value = 61
$ws.Range("AA8").Value = 60.4
$ws.Range("AC8").Value = 57281
$ws.Range("AD8").Value = 567.1
$ws.Range("AE8").Value = 716.3
$ws.Range("AH8").Value = 306
$ws.Range("AI8").Value = 777
$ws.Range("AJ8").Value = 3498
$ws.Range("AL8").Value = 830.2
$ws.Range("AM8").Value = 69
$ws.Range("AN8").Value = 68.3
$ws.Range("AO8").Value = 0.5
# Row 9
$ws.Range("C9").Value = 1929
$ws.Range("D9").Value = 83.90000000000001
$ws.Range("E9").Value = 343.5
$ws.Range("J9").Value = 1648
$ws.Range("L9").Value = 482.2
$ws.Range("P9").Value = 698
$ws.Range("Q9").Value = 30.3
$ws.Range("R9").Value = 133.8
$ws.Range("W9").Value = 643
$ws.Range("Y9").Value = 174.5
$ws.Range("AC9").Value = 2230
$ws.Range("AD9").Value = 97
$ws.Range("AE9").Value = 256.2
$ws.Range("AJ9").Value = 981
$ws.Range("AL9").Value = 446
$ws.Range("AM9").Value = 5
$ws.Range("AN9").Value = 21.7
$ws.Range("AO9").Value = -1.8
# Row 10
$ws.Range("C10").Value = 747611
$ws.Range("D10").Value = 598.6
$ws.Range("E10").Value = 1535.1
$ws.Range("I10").Value = 550
$ws.Range("J10").Value = 24791
$ws.Range("L10").Value = 1139.7
$ws.Range("M10").Value = 656
$ws.Range("N10").Value = 52.5
$ws.Range("P10").Value = 319839
$ws.Range("Q10").Value = 256.1
$ws.Range("R10").Value = 1277.3
$ws.Range("V10").Value = 99
$ws.Range("W10").Value = 24554
$ws.Range("Y10").Value = 473.1
$ws.Range("Z10").Value = 676
$ws.Range("AA10").Value = 54.1
$ws.Range("AC10").Value = 454111
$ws.Range("AD10").Value = 363.6
$ws.Range("AE10").Value = 608.1
$ws.Range("AH10").Value = 59
$ws.Range("AI10").Value = 514
$ws.Range("AJ10").Value = 5226
$ws.Range("AL10").Value = 619.5
$ws.Range("AM10").Value = 733
$ws.Range("AN10").Value = 58.7
$ws.Range("AO10").Value = 0
# Row 11
$ws.Range("C11").Value = 55848
$ws.Range("D11").Value = 3723.2
$ws.Range("E11").Value = 13441.8
$ws.Range("J11").Value = 52263
$ws.Range("L11").Value = 11169.6
$ws.Range("P11").Value = 25527
$ws.Range("Q11").Value = 1701.8
$ws.Range("R11").Value = 6341.5
$ws.Range("V11").Value = 10.5
$ws.Range("W11").Value = 24618
$ws.Range("Y11").Value = 5105.4
$ws.Range("AC11").Value = 1993
$ws.Range("AD11").Value = 132.9
$ws.Range("AE11").Value = 305.3
$ws.Range("AI11").Value = 56
$ws.Range("AJ11").Value = 932
$ws.Range("AL11").Value = 398.6
$ws.Range("AM11").Value = 5
$ws.Range("AN11").Value = 33.3
$ws.Range("AO11").Value = -1.2
# Row 12
$ws.Range("C12").Value = 29308
$ws.Range("D12").Value = 553
$ws.Range("E12").Value = 1332.4
$ws.Range("I12").Value = 368
$ws.Range("J12").Value = 6402
$ws.Range("L12").Value = 1127.2
$ws.Range("O12").Value = 0.3
$ws.Range("P12").Value = 6425
$ws.Range("Q12").Value = 121.2
$ws.Range("R12").Value = 255.3
$ws.Range("V12").Value = 86
$ws.Range("W12").Value = 1182
$ws.Range("Y12").Value = 279.3
$ws.Range("AC12").Value = 20215
$ws.Range("AD12").Value = 381.4
$ws.Range("AE12").Value = 473.1
$ws.Range("AH12").Value = 240
$ws.Range("AI12").Value = 697
$ws.Range("AJ12").Value = 1627
$ws.Range("AL12").Value = 697.1
$ws.Range("AM12").Value = 29
$ws.Range("AN12").Value = 54.7
$ws.Range("AO12").Value = -0.2
# Row 13
$ws.Range("C13").Value = 96445
$ws.Range("D13").Value = 398.5
$ws.Range("E13").Value = 1113.9
$ws.Range("I13").Value = 348
$ws.Range("J13").Value = 9810
$ws.Range("L13").Value = 909.9
$ws.Range("M13").Value = 106
$ws.Range("N13").Value = 43.8
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 33629
$ws.Range("Q13").Value = 139
$ws.Range("R13").Value = 630.2
$ws.Range("V13").Value = 53.8
$ws.Range("W13").Value = 6901
$ws.Range("Y13").Value = 303
$ws.Range("Z13").Value = 111
$ws.Range("AA13").Value = 45.9
$ws.Range("AC13").Value = 73772
$ws.Range("AD13").Value = 304.8
$ws.Range("AE13").Value = 587.6
$ws.Range("AH13").Value = 16
$ws.Range("AI13").Value = 449.2
$ws.Range("AJ13").Value = 6193
$ws.Range("AL13").Value = 567.5
$ws.Range("AM13").Value = 130
$ws.Range("AN13").Value = 53.7
# Row 14
$ws.Range("C14").Value = 42984
$ws.Range("D14").Value = 405.5
$ws.Range("E14").Value = 977.1
$ws.Range("I14").Value = 355.5
$ws.Range("J14").Value = 6478
$ws.Range("L14").Value = 767.6
$ws.Range("M14").Value = 56
$ws.Range("N14").Value = 52.8
$ws.Range("P14").Value = 9034
$ws.Range("Q14").Value = 85.2
$ws.Range("R14").Value = 408.1
$ws.Range("V14").Value = 24.8
$ws.Range("W14").Value = 4083
$ws.Range("Y14").Value = 164.3
$ws.Range("AC14").Value = 34514
$ws.Range("AD14").Value = 325.6
$ws.Range("AE14").Value = 614.5
$ws.Range("AH14").Value = 33
$ws.Range("AI14").Value = 379.8
$ws.Range("AJ14").Value = 3597
$ws.Range("AL14").Value = 585
$ws.Range("AM14").Value = 59
$ws.Range("AN14").Value = 55.7
$ws.Range("AO14").Value = -0.1
